$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Articles")

# Add new header "Fichier" in H1
$ws.Range("H1").Value = 'Fichier'

# Row 2
$ws.Range("C2").Value = 'Dhruv Sarkar'
$ws.Range("E2").Value = '2025-08-23T11:21:24+00:00'
$ws.Range("H2").Value = '2508.16992v1.Online_Learning_for_Approximately_Convex_Functions_with_Long_term_Adversarial_Constraints.pdf'

# Row 3
$ws.Range("C3").Value = 'Anders Aamand'
$ws.Range("E3").Value = '2025-05-29T19:47:09+00:00'
$ws.Range("H3").Value = '2505.23967v1.Improved_Approximations_for_Hard_Graph_Problems_using_Predictions.pdf'

# Row 4
$ws.Range("C4").Value = 'Simone Foderà'
$ws.Range("E4").Value = '2024-09-09T10:07:12+00:00'
$ws.Range("H4").Value = '2409.05475v1.Reinforcement_Learning_for_Variational_Quantum_Circuits_Design.pdf'

# Row 5
$ws.Range("C5").Value = 'Yanhui Zhu'
$ws.Range("E5").Value = '2024-08-08T17:50:16+00:00'
$ws.Range("H5").Value = '2408.04620v2.Regularized_Unconstrained_Weakly_Submodular_Maximization.pdf'

# Row 6
$ws.Range("C6").Value = 'Yang Liu'
$ws.Range("E6").Value = '2024-06-05T22:52:27+00:00'
$ws.Range("H6").Value = '2406.03647v2.Decision_focused_Graph_Neural_Networks_for_Combinatorial_Optimization.pdf'

# Row 7
$ws.Range("C7").Value = 'Morris Yau'
$ws.Range("E7").Value = '2023-10-01T00:12:31+00:00'
$ws.Range("H7").Value = '2310.00526v7.Are_Graph_Neural_Networks_Optimal_Approximation_Algorithms_.pdf'

# Row 8
$ws.Range("C8").Value = 'Davin Choo'
$ws.Range("E8").Value = '2023-01-09T06:25:44+00:00'
$ws.Range("H8").Value = '2301.03180v3.Subset_verification_and_search_algorithms_for_causal_DAGs.pdf'

# Row 9
$ws.Range("E9").Value = '2022-08-12T21:27:20+00:00'
$ws.Range("H9").Value = '2208.06506v3.Optimal_LP_Rounding_and_Linear_Time_Approximation_Algorithms_for_Clustering_Edge_Colored_Hypergraphs.pdf'

# Row 10
$ws.Range("C10").Value = 'David Ireland'
$ws.Range("E10").Value = '2022-05-20T11:54:03+00:00'
$ws.Range("H10").Value = '2205.10106v1.LeNSE__Learning_To_Navigate_Subgraph_Embeddings_for_Large_Scale_Combinatorial_Optimisation.pdf'

# Row 11
$ws.Range("C11").Value = 'Seth Poulsen'
$ws.Range("E11").Value = '2022-04-08T17:44:59+00:00'
$ws.Range("H11").Value = '2204.04196v3.Efficient_Feedback_and_Partial_Credit_Grading_for_Proof_Blocks_Problems.pdf'

# Row 12
$ws.Range("C12").Value = 'Soheil Behnezhad'
$ws.Range("E12").Value = '2021-12-10T09:46:12+00:00'
$ws.Range("H12").Value = '2112.05415v1.Stochastic_Vertex_Cover_with_Few_Queries.pdf'

# Row 13
$ws.Range("C13").Value = 'Martin J. A. Schuetz'
$ws.Range("E13").Value = '2021-07-02T16:54:35+00:00'
$ws.Range("H13").Value = '2107.01188v2.Combinatorial_Optimization_with_Physics_Inspired_Graph_Neural_Networks.pdf'

# Row 14
$ws.Range("C14").Value = 'Lukas Gianinazzi'
$ws.Range("E14").Value = '2021-06-07T13:21:09+00:00'
$ws.Range("H14").Value = '2106.03594v3.Learning_Combinatorial_Node_Labeling_Algorithms.pdf'

# Row 15
$ws.Range("C15").Value = 'Magnús M. Halldórsson'
$ws.Range("E15").Value = '2020-12-17T09:54:24+00:00'
$ws.Range("H15").Value = '2012.09475v2.Query_Competitive_Sorting_with_Uncertainty.pdf'

# Row 16
$ws.Range("C16").Value = 'Faisal N. Abu-Khzam'
$ws.Range("E16").Value = '2020-06-08T15:40:04+00:00'
$ws.Range("H16").Value = '2006.04689v1.Graph_Minors_Meet_Machine_Learning__the_Power_of_Obstructions.pdf'

# Row 17
$ws.Range("C17").Value = 'Yaoxin Li'
$ws.Range("E17").Value = '2020-04-14T14:11:00+00:00'
$ws.Range("H17").Value = '2004.07300v1.Gumbel_softmax_based_Optimization__A_Simple_General_Framework_for_Optimization_Problems_on_Graphs.pdf'

# Row 18
$ws.Range("C18").Value = 'Evripidis Bampis'
$ws.Range("E18").Value = '2019-07-12T20:37:07+00:00'
$ws.Range("H18").Value = '1907.05944v2.Online_learning_for_min_max_discrete_problems.pdf'

# Row 19
$ws.Range("C19").Value = 'Ryoma Sato'
$ws.Range("E19").Value = '2019-05-24T14:41:17+00:00'
$ws.Range("H19").Value = '1905.10261v2.Approximation_Ratios_of_Graph_Neural_Networks_for_Combinatorial_Problems.pdf'

# Row 20
$ws.Range("C20").Value = 'Ceyhun Eksin'
$ws.Range("E20").Value = '2018-12-08T18:51:51+00:00'
$ws.Range("H20").Value = '1812.03366v2.Control_of_learning_in_anti_coordination_network_games.pdf'

# Row 21
$ws.Range("C21").Value = 'Mohsen Ghaffari'
$ws.Range("E21").Value = '2018-07-17T07:01:03+00:00'
$ws.Range("H21").Value = '1807.06251v1.Sparsifying_Distributed_Algorithms_with_Ramifications_in_Massively_Parallel_Computation_and_Centralized_Local_Computation.pdf'

# Row 22
$ws.Range("C22").Value = 'Daniel Selsam'
$ws.Range("E22").Value = '2018-02-11T03:04:28+00:00'
$ws.Range("H22").Value = '1802.03685v4.Learning_a_SAT_Solver_from_Single_Bit_Supervision.pdf'

# Row 23
$ws.Range("C23").Value = 'Hanjun Dai'
$ws.Range("E23").Value = '2017-04-05T23:08:07+00:00'
$ws.Range("H23").Value = '1704.01665v4.Learning_Combinatorial_Optimization_Algorithms_over_Graphs.pdf'
